$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 346.66666
$ws.Range("I4").Value = 264.875
$ws.Range("K4").Value = 264.875
$ws.Range("M4").Value = -150.875

$ws.Range("H17").Value = 2491.1667
$ws.Range("J17").Value = 1996.7
$ws.Range("L17").Value = 5990.1
$ws.Range("N17").Value = -6326.1

$ws.Range("H33").Value = 407.83334
$ws.Range("I33").Value = 199.5
$ws.Range("K33").Value = 199.5
$ws.Range("M33").Value = 29.5

$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 20000
$ws.Range("K44").Value = 20000
$ws.Range("M44").Value = -19538

$ws.Range("H62").Value = 6429.154
$ws.Range("I62").Value = 3916.2
$ws.Range("K62").Value = 3916.2
$ws.Range("M62").Value = -3292.2

$ws.Range("H65").Value = 6429.154
$ws.Range("I65").Value = 3916.2
$ws.Range("K65").Value = 19581
$ws.Range("M65").Value = -16461

$ws.Range("H92").Value = 928
$ws.Range("I92").Value = 920.1111
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 920.1111
$ws.Range("L92").Value = 999
$ws.Range("M92").Value = 327.8889
$ws.Range("N92").Value = -3495

$ws.Range("H96").Value = 3852.375
$ws.Range("J96").Value = 4504.75
$ws.Range("L96").Value = 13514.25
$ws.Range("N96").Value = -16260.25

$ws.Range("H138").Value = 3974.8928
$ws.Range("J138").Value = 4161.481
$ws.Range("L138").Value = 12484.443
$ws.Range("N138").Value = -22764.443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17508.53
$ws.Range("I32").Value = 15895.417
$ws.Range("K32").Value = 15895.417
$ws.Range("M32").Value = -15608.417

$ws.Range("H45").Value = 2035.375
$ws.Range("I45").Value = 1948
$ws.Range("K45").Value = 1948
$ws.Range("M45").Value = -1571

$ws.Range("H61").Value = 1896.75
$ws.Range("I61").Value = 1896.75
$ws.Range("K61").Value = 1896.75
$ws.Range("M61").Value = -1684.75

$ws.Range("H97").Value = 1444.3846
$ws.Range("I97").Value = 305.5
$ws.Range("K97").Value = 305.5
$ws.Range("M97").Value = 190.5

$ws.Range("H102").Value = 1301.4286
$ws.Range("I102").Value = 722
$ws.Range("K102").Value = 722
$ws.Range("M102").Value = 900

$ws.Range("H122").Value = 590156.0600000001
$ws.Range("I122").Value = 668210.2
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 2004630.6
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -2002180.6
$ws.Range("N122").Value = -19150

$ws.Range("H132").Value = 2528.8696
$ws.Range("I132").Value = 2528.8696
$ws.Range("K132").Value = 7586.6088
$ws.Range("M132").Value = -5056.6088

$ws.Range("H136").Value = 1896.75
$ws.Range("I136").Value = 1896.75
$ws.Range("K136").Value = 5690.25
$ws.Range("M136").Value = -3140.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1213.2
$ws.Range("I20").Value = 1247.8334
$ws.Range("K20").Value = 1247.8334
$ws.Range("M20").Value = -1000.8334

$ws.Range("H80").Value = 229.66667
$ws.Range("I80").Value = 322.83334
$ws.Range("K80").Value = 322.83334
$ws.Range("M80").Value = 675.16666

$ws.Range("H83").Value = 229.66667
$ws.Range("I83").Value = 322.83334
$ws.Range("K83").Value = 1614.1667
$ws.Range("M83").Value = 3377.8333

$ws.Range("H94").Value = 2461.2307
$ws.Range("I94").Value = 2381.4546
$ws.Range("J94").Value = 2900
$ws.Range("K94").Value = 2381.4546
$ws.Range("L94").Value = 2900
$ws.Range("M94").Value = -1930.4546
$ws.Range("N94").Value = -3802

$ws.Range("H105").Value = 5255.625
$ws.Range("I105").Value = 3774.3572
$ws.Range("J105").Value = 15624.5
$ws.Range("K105").Value = 3774.3572
$ws.Range("L105").Value = 15624.5
$ws.Range("M105").Value = -2027.3572
$ws.Range("N105").Value = -19118.5

$ws.Range("H107").Value = 1034.5
$ws.Range("I107").Value = 1034.5
$ws.Range("K107").Value = 1034.5
$ws.Range("M107").Value = 885.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 163.66667
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H62").Value = 102249
$ws.Range("J62").Value = 135332.33
$ws.Range("L62").Value = 135332.33
$ws.Range("N62").Value = -136580.33

$ws.Range("H65").Value = 102249
$ws.Range("J65").Value = 135332.33
$ws.Range("L65").Value = 676661.6499999999
$ws.Range("N65").Value = -682901.6499999999

$ws.Range("H105").Value = 2041.5454
$ws.Range("J105").Value = 3099.8
$ws.Range("L105").Value = 3099.8
$ws.Range("N105").Value = -6593.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 584.5
$ws.Range("I11").Value = 169
$ws.Range("K11").Value = 507
$ws.Range("M11").Value = -367

$ws.Range("H64").Value = 41667228
$ws.Range("I64").Value = 495.4
$ws.Range("J64").Value = 111111780
$ws.Range("K64").Value = 1486.2
$ws.Range("L64").Value = 333335340
$ws.Range("M64").Value = -1216.2
$ws.Range("N64").Value = -333335880

$ws.Range("H67").Value = 41667228
$ws.Range("I67").Value = 495.4
$ws.Range("J67").Value = 111111780
$ws.Range("K67").Value = 1486.2
$ws.Range("L67").Value = 333335340
$ws.Range("M67").Value = -550.1999999999998
$ws.Range("N67").Value = -333337212

$ws.Range("H75").Value = 362
$ws.Range("I75").Value = 386.5
$ws.Range("K75").Value = 1159.5
$ws.Range("M75").Value = -161.5

$ws.Range("H78").Value = 362
$ws.Range("I78").Value = 386.5
$ws.Range("K78").Value = 3478.5
$ws.Range("M78").Value = 1513.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6002.5
$ws.Range("I70").Value = 6002.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6002.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5732.5
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 6002.5
$ws.Range("I73").Value = 6002.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6002.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5066.5
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -84984

$ws.Range("H97").Value = 1266.3125
$ws.Range("I97").Value = 1013.3333
$ws.Range("J97").Value = 1591.5714
$ws.Range("K97").Value = 1013.3333
$ws.Range("L97").Value = 1591.5714
$ws.Range("M97").Value = -517.3333
$ws.Range("N97").Value = -2583.5714

$ws.Range("H113").Value = 4656.091
$ws.Range("I113").Value = 3840.3333
$ws.Range("J113").Value = 4962
$ws.Range("K113").Value = 3840.3333
$ws.Range("L113").Value = 4962
$ws.Range("M113").Value = -1670.3333
$ws.Range("N113").Value = -9302

$ws.Range("H122").Value = 169463.17
$ws.Range("I122").Value = 2250.75
$ws.Range("J122").Value = 503888
$ws.Range("K122").Value = 6752.25
$ws.Range("L122").Value = 1511664
$ws.Range("M122").Value = -4302.25
$ws.Range("N122").Value = -1516564

$ws.Range("H132").Value = 1927.7188
$ws.Range("I132").Value = 1514.3214
$ws.Range("K132").Value = 4542.9642
$ws.Range("M132").Value = -2012.9642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 582.9167
$ws.Range("I93").Value = 443.8889
$ws.Range("K93").Value = 443.8889
$ws.Range("M93").Value = 804.1111000000001

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 790.63635
$ws.Range("I81").Value = 834.7
$ws.Range("K81").Value = 1669.4
$ws.Range("M81").Value = -608.4000000000001

$ws.Range("H84").Value = 790.63635
$ws.Range("I84").Value = 834.7
$ws.Range("K84").Value = 8347
$ws.Range("M84").Value = -3043

$ws.Range("H96").Value = 1086.5
$ws.Range("I96").Value = 1033.909
$ws.Range("K96").Value = 1033.909
$ws.Range("M96").Value = 339.0909999999999

$ws.Range("H122").Value = 16329.333
$ws.Range("I122").Value = 16329.333
$ws.Range("K122").Value = 48987.999
$ws.Range("M122").Value = -46537.999
